$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Drill -----------------------------------------------------
# Set values in the same order the strings were first introduced so that
# the shared-strings table is rebuilt in the same order as the target file
# (Drill, Utility/Movement, description, comment, Boxing Glove, description).
$ws.Range("A3").Value = "Drill"
$ws.Range("C3").Value = "Utility/Movement"
$ws.Range("B3").Value = "When activated, the tank using it moves in a user-specified direction, destroying terrain in its path for a set time. Gravity affects the tank as usual while the item is activated, to prevent tanks from using the drill to fly."
$ws.Range("F3").Value = "Ben: Feels like a given once we get destructable terrain working. Sufficiently unique experience to be meet complexity requirement."

# Approval Status cell needs the same "To be reviewed" highlight style as D2.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "To be reviewed"

$ws.Range("E3").Value = "Medium"

# --- Row 4: Boxing Glove ----------------------------------------------
$ws.Range("A4").Value = "Boxing Glove"
$ws.Range("B4").Value = "A Melee weapon that has short range and low damage, but high knockback."
$ws.Range("C4").Value = "Weapon"

$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "To be reviewed"

$ws.Range("E4").Value = "Medium"

$excel.CutCopyMode = 0

# --- Misc UI state ------------------------------------------------------
[void]$ws.Range("F4").Select()

Write-Output "done"
